$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.019999999999999
$ws.Cells.Item(2, 3).Value = 1.173724721177118
$ws.Cells.Item(2, 4).Value = 1.180143870303751
$ws.Cells.Item(2, 5).Value = 1.165107498415543
$ws.Cells.Item(2, 6).Value = 1.177269549995638
$ws.Cells.Item(2, 9).Value = 1.034750151734215
$ws.Cells.Item(2, 10).Value = 1.178182774779943
$ws.Cells.Item(2, 11).Value = 1.182580264568403
$ws.Cells.Item(2, 12).Value = 1.167576858788469
$ws.Cells.Item(2, 13).Value = 1.179712176442149
$ws.Cells.Item(2, 14).Value = 1.179855929459584

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.178499317574796
$ws.Cells.Item(3, 4).Value = 1.18479146393432
$ws.Cells.Item(3, 5).Value = 1.169412856507996
$ws.Cells.Item(3, 6).Value = 1.181705452232316
$ws.Cells.Item(3, 9).Value = 1.034968540629482
$ws.Cells.Item(3, 10).Value = 1.182635718723106
$ws.Cells.Item(3, 11).Value = 1.187055390228695
$ws.Cells.Item(3, 12).Value = 1.171707895353263
$ws.Cells.Item(3, 13).Value = 1.183975552772438
$ws.Cells.Item(3, 14).Value = 1.184315197093905

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.181555438499013
$ws.Cells.Item(4, 4).Value = 1.187765789567637
$ws.Cells.Item(4, 5).Value = 1.172167173474739
$ws.Cells.Item(4, 6).Value = 1.184543082489814
$ws.Cells.Item(4, 9).Value = 1.035104774096685
$ws.Cells.Item(4, 10).Value = 1.185484316741865
$ws.Cells.Item(4, 11).Value = 1.189918173745531
$ws.Cells.Item(4, 12).Value = 1.174349426641909
$ws.Cells.Item(4, 13).Value = 1.186701569276752
$ws.Cells.Item(4, 14).Value = 1.187167840448591

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.18283253482353
$ws.Cells.Item(5, 4).Value = 1.189008583133323
$ws.Cells.Item(5, 5).Value = 1.173317802325079
$ws.Cells.Item(5, 6).Value = 1.185728466275991
$ws.Cells.Item(5, 9).Value = 1.035160849307454
$ws.Cells.Item(5, 10).Value = 1.186674300355982
$ws.Cells.Item(5, 11).Value = 1.191114081693979
$ws.Cells.Item(5, 12).Value = 1.175452636879291
$ws.Cells.Item(5, 13).Value = 1.187840029446056
$ws.Cells.Item(5, 14).Value = 1.188359513975932

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.183046521816671
$ws.Cells.Item(6, 4).Value = 1.189216815241053
$ws.Cells.Item(6, 5).Value = 1.173510578231852
$ws.Cells.Item(6, 6).Value = 1.185927062061119
$ws.Cells.Item(6, 9).Value = 1.035170194895599
$ws.Cells.Item(6, 10).Value = 1.18687366791698
$ws.Cells.Item(6, 11).Value = 1.191314441699809
$ws.Cells.Item(6, 12).Value = 1.175637450604206
$ws.Cells.Item(6, 13).Value = 1.188030746294929
$ws.Cells.Item(6, 14).Value = 1.188559164661735

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.181572532988299
$ws.Cells.Item(7, 4).Value = 1.187782425382114
$ws.Cells.Item(7, 5).Value = 1.172182576518563
$ws.Cells.Item(7, 6).Value = 1.184558950982279
$ws.Cells.Item(7, 9).Value = 1.035105528056915
$ws.Cells.Item(7, 10).Value = 1.185500246736628
$ws.Cells.Item(7, 11).Value = 1.189934183054627
$ws.Cells.Item(7, 12).Value = 1.174364196098042
$ws.Cells.Item(7, 13).Value = 1.186716810786684
$ws.Cells.Item(7, 14).Value = 1.187183793065773

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.175345404105126
$ws.Cells.Item(8, 4).Value = 1.181721547670009
$ws.Cells.Item(8, 5).Value = 1.166569206382772
$ws.Cells.Item(8, 6).Value = 1.178775620011999
$ws.Cells.Item(8, 9).Value = 1.034825020247594
$ws.Cells.Item(8, 10).Value = 1.179694618371821
$ws.Cells.Item(8, 11).Value = 1.184099639798525
$ws.Cells.Item(8, 12).Value = 1.168979647581369
$ws.Cells.Item(8, 13).Value = 1.181159931992742
$ws.Cells.Item(8, 14).Value = 1.181369920042774

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.164103367218748
$ws.Cells.Item(9, 4).Value = 1.17077577832289
$ws.Cells.Item(9, 5).Value = 1.156423872949547
$ws.Cells.Item(9, 6).Value = 1.168321595870597
$ws.Cells.Item(9, 9).Value = 1.034290965616229
$ws.Cells.Item(9, 10).Value = 1.16920069947522
$ws.Cells.Item(9, 11).Value = 1.173553445665439
$ws.Cells.Item(9, 12).Value = 1.15923791598218
$ws.Cells.Item(9, 13).Value = 1.171105415230106
$ws.Cells.Item(9, 14).Value = 1.170861098577669

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.156408105685987
$ws.Cells.Item(10, 4).Value = 1.163280795224431
$ws.Cells.Item(10, 5).Value = 1.149471712575498
$ws.Cells.Item(10, 6).Value = 1.161156972584162
$ws.Cells.Item(10, 9).Value = 1.033906927543563
$ws.Cells.Item(10, 10).Value = 1.162008855881951
$ws.Cells.Item(10, 11).Value = 1.166325803633457
$ws.Cells.Item(10, 12).Value = 1.152555515735182
$ws.Cells.Item(10, 13).Value = 1.164207882478372
$ws.Cells.Item(10, 14).Value = 1.163659041741582

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.153023775512674
$ws.Cells.Item(11, 4).Value = 1.159983971162742
$ws.Cells.Item(11, 5).Value = 1.146412404668183
$ws.Cells.Item(11, 6).Value = 1.158003970033121
$ws.Cells.Item(11, 9).Value = 1.033733694184499
$ws.Cells.Item(11, 10).Value = 1.158843867306624
$ws.Cells.Item(11, 11).Value = 1.163145075327977
$ws.Cells.Item(11, 12).Value = 1.149613281951436
$ws.Cells.Item(11, 13).Value = 1.161170813674561
$ws.Cells.Item(11, 14).Value = 1.160489558519449

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.151758422234229
$ws.Cells.Item(12, 4).Value = 1.158751250795198
$ws.Cells.Item(12, 5).Value = 1.145268305199419
$ws.Cells.Item(12, 6).Value = 1.156824803024507
$ws.Cells.Item(12, 9).Value = 1.033668277609773
$ws.Cells.Item(12, 10).Value = 1.157660213565912
$ws.Cells.Item(12, 11).Value = 1.161955538050469
$ws.Cells.Item(12, 12).Value = 1.148512718039964
$ws.Cells.Item(12, 13).Value = 1.160034761097571
$ws.Cells.Item(12, 14).Value = 1.159304223854658

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.15203022559354
$ws.Cells.Item(13, 4).Value = 1.159016048245768
$ws.Cells.Item(13, 5).Value = 1.145514074782588
$ws.Cells.Item(13, 6).Value = 1.157078106884141
$ws.Cells.Item(13, 9).Value = 1.033682358555131
$ws.Cells.Item(13, 10).Value = 1.157914481487043
$ws.Cells.Item(13, 11).Value = 1.162211069709428
$ws.Cells.Item(13, 12).Value = 1.148749146730608
$ws.Cells.Item(13, 13).Value = 1.160278814284192
$ws.Cells.Item(13, 14).Value = 1.159558852865402

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.152919352039043
$ws.Cells.Item(14, 4).Value = 1.159882242431018
$ws.Cells.Item(14, 5).Value = 1.146317993119491
$ws.Cells.Item(14, 6).Value = 1.157906665261256
$ws.Cells.Item(14, 9).Value = 1.033728308818411
$ws.Cells.Item(14, 10).Value = 1.158746192393305
$ws.Cells.Item(14, 11).Value = 1.163046914843803
$ws.Cells.Item(14, 12).Value = 1.149522468011917
$ws.Cells.Item(14, 13).Value = 1.161077071650566
$ws.Cells.Item(14, 14).Value = 1.16039174489655

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.15346606457716
$ws.Cells.Item(15, 4).Value = 1.160414843124553
$ws.Cells.Item(15, 5).Value = 1.146812276980308
$ws.Cells.Item(15, 6).Value = 1.158416095210268
$ws.Cells.Item(15, 9).Value = 1.033756477677535
$ws.Cells.Item(15, 10).Value = 1.159257559993685
$ws.Cells.Item(15, 11).Value = 1.163560824724753
$ws.Cells.Item(15, 12).Value = 1.149997906789024
$ws.Cells.Item(15, 13).Value = 1.161567839248692
$ws.Cells.Item(15, 14).Value = 1.160903838697577

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.156631577106379
$ws.Cells.Item(16, 4).Value = 1.163498476636917
$ws.Cells.Item(16, 5).Value = 1.149673685033457
$ws.Cells.Item(16, 6).Value = 1.161365126633817
$ws.Cells.Item(16, 9).Value = 1.033918275931879
$ws.Cells.Item(16, 10).Value = 1.162217800813323
$ws.Cells.Item(16, 11).Value = 1.166535788024768
$ws.Cells.Item(16, 12).Value = 1.152749724737988
$ws.Cells.Item(16, 13).Value = 1.164408349094215
$ws.Cells.Item(16, 14).Value = 1.163868283398722

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.158602948967022
$ws.Cells.Item(17, 4).Value = 1.165418705795594
$ws.Cells.Item(17, 5).Value = 1.151455197145776
$ws.Cells.Item(17, 6).Value = 1.163201141209703
$ws.Cells.Item(17, 9).Value = 1.034017890418512
$ws.Cells.Item(17, 10).Value = 1.164060790161449
$ws.Cells.Item(17, 11).Value = 1.168387947621326
$ws.Cells.Item(17, 12).Value = 1.154462571378658
$ws.Cells.Item(17, 13).Value = 1.166176372709852
$ws.Cells.Item(17, 14).Value = 1.165713890003117

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.159747806362424
$ws.Cells.Item(18, 4).Value = 1.166533807383404
$ws.Cells.Item(18, 5).Value = 1.152489623978858
$ws.Cells.Item(18, 6).Value = 1.16426719540942
$ws.Cells.Item(18, 9).Value = 1.034075325921422
$ws.Cells.Item(18, 10).Value = 1.165130893708505
$ws.Cells.Item(18, 11).Value = 1.169463377067138
$ws.Cells.Item(18, 12).Value = 1.155456971857996
$ws.Cells.Item(18, 13).Value = 1.167202795481157
$ws.Cells.Item(18, 14).Value = 1.166785513219953

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.160137335561016
$ws.Cells.Item(19, 4).Value = 1.166913202981694
$ws.Cells.Item(19, 5).Value = 1.152841550934811
$ws.Cells.Item(19, 6).Value = 1.164629879135741
$ws.Cells.Item(19, 9).Value = 1.034094797445955
$ws.Cells.Item(19, 10).Value = 1.165494955084952
$ws.Cells.Item(19, 11).Value = 1.169829250580711
$ws.Cells.Item(19, 12).Value = 1.155795254844466
$ws.Cells.Item(19, 13).Value = 1.16755197014721
$ws.Cells.Item(19, 14).Value = 1.167150091605313

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.15839196061445
$ws.Cells.Item(20, 4).Value = 1.165213196746586
$ws.Cells.Item(20, 5).Value = 1.151264546473827
$ws.Cells.Item(20, 6).Value = 1.163004659871394
$ws.Cells.Item(20, 9).Value = 1.03400727200591
$ws.Cells.Item(20, 10).Value = 1.163863562485613
$ws.Cells.Item(20, 11).Value = 1.168189738447351
$ws.Cells.Item(20, 12).Value = 1.154279285206679
$ws.Cells.Item(20, 13).Value = 1.165987183347414
$ws.Cells.Item(20, 14).Value = 1.165516382241359

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.152657757982279
$ws.Cells.Item(21, 4).Value = 1.159627397690814
$ws.Cells.Item(21, 5).Value = 1.146081475897408
$ws.Cells.Item(21, 6).Value = 1.157662899566575
$ws.Cells.Item(21, 9).Value = 1.033714807371855
$ws.Cells.Item(21, 10).Value = 1.158501499350584
$ws.Cells.Item(21, 11).Value = 1.162801005408242
$ws.Cells.Item(21, 12).Value = 1.149294959432404
$ws.Cells.Item(21, 13).Value = 1.160842227351687
$ws.Cells.Item(21, 14).Value = 1.160146704361642

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.149004412845144
$ws.Cells.Item(22, 4).Value = 1.156068113637364
$ws.Cells.Item(22, 5).Value = 1.142777712791172
$ws.Cells.Item(22, 6).Value = 1.154257822241259
$ws.Cells.Item(22, 9).Value = 1.033524718223271
$ws.Cells.Item(22, 10).Value = 1.155083454233854
$ws.Cells.Item(22, 11).Value = 1.159365976563894
$ws.Cells.Item(22, 12).Value = 1.146116443445427
$ws.Cells.Item(22, 13).Value = 1.157561191383297
$ws.Cells.Item(22, 14).Value = 1.15672380522879

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.150945816481299
$ws.Cells.Item(23, 4).Value = 1.157959578284751
$ws.Cells.Item(23, 5).Value = 1.14453349311334
$ws.Cells.Item(23, 6).Value = 1.15606746065074
$ws.Cells.Item(23, 9).Value = 1.0336260856173
$ws.Cells.Item(23, 10).Value = 1.156899987928902
$ws.Cells.Item(23, 11).Value = 1.161191534487295
$ws.Cells.Item(23, 12).Value = 1.147805797438571
$ws.Cells.Item(23, 13).Value = 1.159305041263326
$ws.Cells.Item(23, 14).Value = 1.158542918610045

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.158487312623846
$ws.Cells.Item(24, 4).Value = 1.165306072664002
$ws.Cells.Item(24, 5).Value = 1.151350707800977
$ws.Cells.Item(24, 6).Value = 1.163093456328698
$ws.Cells.Item(24, 9).Value = 1.034012072071703
$ws.Cells.Item(24, 10).Value = 1.163952696238842
$ws.Cells.Item(24, 11).Value = 1.168279315768188
$ws.Cells.Item(24, 12).Value = 1.15436211875646
$ws.Cells.Item(24, 13).Value = 1.166072684793261
$ws.Cells.Item(24, 14).Value = 1.165605642574741

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.16704355719825
$ws.Cells.Item(25, 4).Value = 1.173638929531014
$ws.Cells.Item(25, 5).Value = 1.159078563614006
$ws.Cells.Item(25, 6).Value = 1.171057229047690
$ws.Cells.Item(25, 9).Value = 1.034433870437329
$ws.Cells.Item(25, 10).Value = 1.171946739734431
$ws.Cells.Item(25, 11).Value = 1.176313163029073
$ws.Cells.Item(25, 12).Value = 1.161788176606312
$ws.Cells.Item(25, 13).Value = 1.17373766673407
$ws.Cells.Item(25, 14).Value = 1.173611038529024
